# Apply the "fixed save origin test" edit: update the computed test-state
# values (columns B-F) on rows 3,4,7,8,11,12,15,16 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 5
$ws.Range("D3").Value = 176.3651844588868
$ws.Range("E3").Value = 0.8
$ws.Range("F3").Value = 0.008352056178274002

# Row 4
$ws.Range("D4").Value = 378.9637784732679
$ws.Range("E4").Value = 0.8
$ws.Range("F4").Value = 0.008352056178274002

# Row 7
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 328.1082605273225
$ws.Range("E7").Value = 0.8
$ws.Range("F7").Value = 0.008352056178274002

# Row 8
$ws.Range("D8").Value = 714.3960623361766
$ws.Range("E8").Value = 0.8
$ws.Range("F8").Value = 0.008352056178274002

# Row 11
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 481.9949231290301
$ws.Range("E11").Value = 0.8
$ws.Range("F11").Value = 0.008352056178274002

# Row 12
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 1053.149130220529
$ws.Range("E12").Value = 0.8
$ws.Range("F12").Value = 0.008352056178274002

# Row 15
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 559.0635379518063
$ws.Range("E15").Value = 0.8
$ws.Range("F15").Value = 0.008352056178274002

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1222.999659563553
$ws.Range("E16").Value = 0.8
$ws.Range("F16").Value = 0.008352056178274002
